$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22 previously had no explicit cell style; bring it in line with the
#     rest of the data rows by copying the formatting already used by row 21.
$ws.Range("A22:F22").Style = $ws.Range("A21:F21").Style

# --- Append the new sale recorded in row 23.
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "2026-01-21 13:31:40"
$ws.Range("C23").Value = "A001"
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 154.25
$ws.Range("F23").Value = 308.5

# Column B carries a column-level style (used for the date column), which
# Excel would otherwise apply automatically to the freshly written B23 cell.
# The new row is unformatted (just like row 22 was before this edit), so
# reset B23 back to the plain/default style used by the rest of row 23.
$ws.Range("B23").Style = $ws.Range("A23").Style
